$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new incoming-mail row (row 21) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A21").Value = "Status van mijn bestelling"
$logs.Range("B21").Value = "mailmind.test@zohomail.eu"
$logs.Range("C21").Value = "Wanneer wordt mijn bestelling bezorgd?"
$logs.Range("D21").Value = "Bestelling / Levering"
$logs.Range("E21").Value = "Beste klant,`nHartelijk dank voor uw bericht. Om u te kunnen helpen met het checken van de bezorgdatum van uw bestelling, hebben we wat meer informatie nodig. Kunt u ons alstublieft uw bestelnummer doorgeven? Met deze informatie kunnen wij direct voor u nakijken wanneer uw bestelling wordt bezorgd.`nAlvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Naam] `nKlantenservice Team"
$logs.Range("F21").Value = "2025-06-24 20:16:49"
$logs.Range("G21").Value = "Ja"

# Extend the conditional-formatting ranges to cover the new row
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))

# --- Dashboard sheet: re-sort the category counts now that
# "Bestelling / Levering" has grown from 2 to 3 mails ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Bestelling / Levering"
$dash.Range("B5").Value = 3

$dash.Range("A6").Value = "Sollicitatie / Vacature"
$dash.Range("B6").Value = 2

$dash.Range("A8").Value = "Offerte / Prijsaanvraag"
$dash.Range("B8").Value = 2
